{"js": "// Highlight quantitative impact metrics (percentages, dollar amounts, large\n// numbers) across the resume bullets by splitting the affected runs and\n// applying bold + a dark slate color (#2C3E50) to just the metric text.\n\nconst METRIC_COLOR = \"#2C3E50\";\nconst BULLET = \"\\u2022\"; // \u2022\nconst PM = \"\\u00B1\"; // \u00b1\n\n// Ordered list of {exactText, metrics[]} describing each paragraph (matched\n// by its full, exact text) and the ordered substrings inside it that should\n// become bold + colored.\nconst edits = [\n  {\n    exactText:\n      BULLET +\n      \" Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    metrics: [\"23%\", \"64%\"],\n  },\n  {\n    exactText:\n      BULLET +\n      \" Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \" +\n      PM +\n      \"4.2% to \" +\n      PM +\n      \"2.1%\",\n    metrics: [\"87%\", \"71%\", PM + \"4.2%\", PM + \"2.1%\"],\n  },\n  {\n    exactText:\n      BULLET +\n      \" Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    metrics: [\"1,200\"],\n  },\n  {\n    exactText:\n      BULLET +\n      \" Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    metrics: [\"$400M\", \"$1B\"],\n  },\n  {\n    exactText:\n      BULLET +\n      \" Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    metrics: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    exactText:\n      BULLET +\n      \" Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    metrics: [\"87%\", \"71%\"],\n  },\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const edit of edits) {\n  // Find the (unique) paragraph whose full text exactly matches.\n  let target = null;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text === edit.exactText) {\n      target = paragraphs.items[i];\n      break;\n    }\n  }\n  if (!target) {\n    continue;\n  }\n\n  // Apply bold + color to each metric substring, in order, so repeated\n  // numbers (e.g. \"23%\" then \"64%\") each match their own occurrence.\n  for (const metric of edit.metrics) {\n    const results = target.search(metric, { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n\n    if (results.items.length > 0) {\n      const hit = results.items[0];\n      hit.font.bold = true;\n      hit.font.color = METRIC_COLOR;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Highlight quantitative impact metrics (percentages, dollar amounts, large\n# numbers) across the resume bullets by splitting the affected runs and\n# applying bold + a dark slate color (#2C3E50) to just the metric text.\n\n$d = $word.ActiveDocument\n\n# RGB(0x2C,0x3E,0x50) expressed as the BGR integer Word's Font.Color expects.\n$metricColor = 5258796\n$bullet = [char]0x2022\n$pm = [char]0xB1\n$cr = [char]0x0D\n\nfunction Set-MetricHighlights {\n    param(\n        [string]$ExactText,   # full paragraph text (no leading bullet trimming, no trailing CR)\n        [string[]]$Metrics    # ordered list of substrings (within that paragraph) to bold+color\n    )\n    # NOTE: always call this positionally (Set-MetricHighlights $text $metrics) --\n    # named parameter binding is not reliable in this host.\n\n    $hit = $null\n    foreach ($para in $d.Paragraphs) {\n        $t = $para.Range.Text.TrimEnd($cr)\n        if ($t -eq $ExactText) {\n            $hit = $para\n            break\n        }\n    }\n    if ($null -eq $hit) {\n        return\n    }\n\n    $paraRange = $hit.Range\n    $searchStart = $paraRange.Start\n    $paraEnd = $paraRange.End\n\n    foreach ($metric in $Metrics) {\n        $r = $d.Range($searchStart, $paraEnd)\n        $found = $r.Find.Execute($metric)\n        if ($found) {\n            $r.Font.Bold = 1\n            $r.Font.Color = $metricColor\n            $searchStart = $r.End\n        }\n    }\n}\n\nSet-MetricHighlights ($bullet + \" Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\") @(\"23%\", \"64%\")\n\nSet-MetricHighlights ($bullet + \" Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \" + $pm + \"4.2% to \" + $pm + \"2.1%\") @(\"87%\", \"71%\", ($pm + \"4.2%\"), ($pm + \"2.1%\"))\n\nSet-MetricHighlights ($bullet + \" Wrote RFP and analyzed bids from 1,200 vendors for research platform development\") @(\"1,200\")\n\nSet-MetricHighlights ($bullet + ' Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+') @('$400M', '$1B')\n\nSet-MetricHighlights ($bullet + ' Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M') @(\"73.5%\", '$4.7M')\n\nSet-MetricHighlights ($bullet + \" Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\") @(\"87%\", \"71%\")\n"}
